$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = '27.857.89'
$ws.Range('E2').Value = '  -0.42%  '

# Row 3
$ws.Range('D3').Value = '1.908.01'
$ws.Range('E3').Value = '  +0.16%  '

# Row 4
$ws.Range('D4').Value = '1.0000'
$ws.Range('E4').Value = '  -0.38%  '

# Row 5
$ws.Range('D5').Value = '312.65'
$ws.Range('E5').Value = '  -1.48%  '

# Row 6
$ws.Range('D6').Value = '0.9991'

# Row 7
$ws.Range('D7').Value = '0.4995'
$ws.Range('E7').Value = '  +3.71%  '

# Row 8
$ws.Range('D8').Value = '0.3796'
$ws.Range('E8').Value = '  -0.18%  '

# Row 9
$ws.Range('D9').Value = '0.07277'
$ws.Range('E9').Value = '  -1.07%  '

# Row 10
$ws.Range('D10').Value = '21.25'
$ws.Range('E10').Value = '  +2.36%  '

# Row 11
$ws.Range('D11').Value = '0.9090'
$ws.Range('E11').Value = '  -2.35%  '

# Row 12
$ws.Range('B12').Value = 'WrappedEther'
$ws.Range('C12').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D12').Value = '1.932.60'
$ws.Range('E12').Value = '  +2.69%  '

# Row 13
$ws.Range('B13').Value = 'TRON'
$ws.Range('C13').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range('D13').Value = '0.07640'
$ws.Range('E13').Value = '  -1.35%  '

# Row 14
$ws.Range('D14').Value = '5.467'
$ws.Range('E14').Value = '  -0.20%  '

# Row 15
$ws.Range('D15').Value = '92.35'
$ws.Range('E15').Value = '  +0.82%  '

# Row 16
$ws.Range('D16').Value = '1.000'
$ws.Range('E16').Value = '  -0.34%  '

# Row 17
$ws.Range('D17').Value = '0.000008719'
$ws.Range('E17').Value = '  -1.76%  '

# Row 18
$ws.Range('D18').Value = '0.9985'
$ws.Range('E18').Value = '  -0.39%  '

# Row 19
$ws.Range('D19').Value = '27.896.44'
$ws.Range('E19').Value = '  -0.40%  '

# Row 20
$ws.Range('D20').Value = '14.65'
$ws.Range('E20').Value = '  -0.14%  '

# Row 21
$ws.Range('D21').Value = '5.168'
$ws.Range('E21').Value = '  +0.57%  '

# Row 22
$ws.Range('D22').Value = '2.164.43'
$ws.Range('E22').Value = '  +0.91%  '

# Row 23
$ws.Range('D23').Value = '10.87'
$ws.Range('E23').Value = '  -0.22%  '

# Row 24
$ws.Range('D24').Value = '6.596'
$ws.Range('E24').Value = '  -0.58%  '

# Row 25
$ws.Range('D25').Value = '152.57'
$ws.Range('E25').Value = '  -2.13%  '

# Row 26
$ws.Range('D26').Value = '1.843'
$ws.Range('E26').Value = '  -3.22%  '

# Row 27
$ws.Range('D27').Value = '2.218'
$ws.Range('E27').Value = '  +4.91%  '

# Row 28
$ws.Range('D28').Value = '18.38'
$ws.Range('E28').Value = '  -0.54%  '

# Row 29
$ws.Range('E29').Value = '  -1.85%  '

# Row 30
$ws.Range('D30').Value = '4.883'
$ws.Range('E30').Value = '  -1.71%  '

# Row 31
$ws.Range('D31').Value = '0.08967'
$ws.Range('E31').Value = '  +0.32%  '

# Row 32
$ws.Range('D32').Value = '3.160'
$ws.Range('E32').Value = '  -2.54%  '

# Row 33
$ws.Range('B33').Value = 'ARBITRUM'
$ws.Range('C33').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D33').Value = '1.237'
$ws.Range('E33').Value = '  -1.01%  '

# Row 34
$ws.Range('B34').Value = 'Filecoin'
$ws.Range('C34').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D34').Value = '4.816'
$ws.Range('E34').Value = '  +3.27%  '

# Row 35
$ws.Range('E35').Value = '  +2.32%  '

# Row 36
$ws.Range('D36').Value = '2.655'
$ws.Range('E36').Value = '  +2.21%  '

# Row 37
$ws.Range('D37').Value = '0.02079'
$ws.Range('E37').Value = '  +1.47%  '

# Row 38
$ws.Range('D38').Value = '3.063'
$ws.Range('E38').Value = '  +2.28%  '

# Row 39
$ws.Range('E39').Value = '  -1.28%  '

# Row 40
$ws.Range('D40').Value = '0.5525'
$ws.Range('E40').Value = '  +0.33%  '

# Row 41
$ws.Range('D41').Value = '0.05280'
$ws.Range('E41').Value = '  -0.03%  '

# Row 42
$ws.Range('D42').Value = '6.800'
$ws.Range('E42').Value = '  -2.56%  '

# Row 43
$ws.Range('D43').Value = '113.69'
$ws.Range('E43').Value = '  +3.20%  '

# Row 44
$ws.Range('D44').Value = '8.510'
$ws.Range('E44').Value = '  +0.38%  '

# Row 45
$ws.Range('D45').Value = '0.1510'
$ws.Range('E45').Value = '  -1.14%  '

# Row 46
$ws.Range('D46').Value = '10.59'
$ws.Range('E46').Value = '  -1.11%  '

# Row 47
$ws.Range('D47').Value = '0.4818'
$ws.Range('E47').Value = '  -0.03%  '

# Row 48
$ws.Range('D48').Value = '0.9986'
$ws.Range('E48').Value = '  -0.39%  '

# Row 49
$ws.Range('D49').Value = '1.635'
$ws.Range('E49').Value = '  -0.31%  '

# Row 50
$ws.Range('D50').Value = '67.32'
$ws.Range('E50').Value = '  -0.60%  '

# Row 51
$ws.Range('D51').Value = '0.06040'
